# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Nueva Caledonia" / "Belice" rows (19xdata stays with the
#     country it describes, so names AND their stats trade places) ---
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("H192").Value = 2

$ws.Range("A193").Value = "Nueva Caledonia"
$ws.Range("D193").Value = 18
$ws.Range("H193").Value = 0

# --- Swap the "Islas Virgenes Britanicas" / "Butan" rows ---
$ws.Range("A212").Value = "Butan"
$ws.Range("D212").Value = 5
$ws.Range("H212").Value = 0

$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 4
$ws.Range("H213").Value = 1

# --- Updated statistics for Arabia Saudita (row 20) ---
$ws.Range("B20").Value = 39048
$ws.Range("C20").Value = 1912
$ws.Range("D20").Value = 11457
$ws.Range("E20").Value = 27345
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 246

# --- Updated statistics for Portugal (row 25) ---
$ws.Range("B25").Value = 27581
$ws.Range("C25").Value = 175
$ws.Range("D25").Value = 2549
$ws.Range("E25").Value = 23897
$ws.Range("F25").Value = 112
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 1135

# --- Updated statistics for Republica de Macedonia (row 86) ---
$ws.Range("B86").Value = 1642
$ws.Range("C86").Value = 20
$ws.Range("D86").Value = 1136
$ws.Range("E86").Value = 415

# --- Updated statistics for Sri Lanka (row 104) ---
$ws.Range("B104").Value = 855
$ws.Range("C104").Value = 8
$ws.Range("E104").Value = 525
